# Sprint 2 solution 1.0
#
# Replaces the two "Sprint 2" bullet paragraphs:
#   "Create the GameObject abstract class."
#   "Create the Jewel class and inherit from the GameObject class."
# with six bullet paragraphs describing Object/Jewel/Furniture/GallifrianMirror,
# a blank bookmarked paragraph, and a stray "ff" paragraph - matching the
# authored Sprint-2 task list.

$d = $word.ActiveDocument

# --- Locate the two source paragraphs robustly (by content, not a hard index) ---
$firstIdx = 0
$secondIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Create the GameObject abstract class.*") {
        $firstIdx = $i
    }
    if ($t -like "Create the Jewel class and inherit from the GameObject class.*") {
        $secondIdx = $i
    }
}

if ($firstIdx -eq 0 -or $secondIdx -eq 0) {
    Write-Host "Could not locate Sprint 2 source paragraphs (first=$firstIdx second=$secondIdx)"
} else {
    $p1 = $d.Paragraphs.Item($firstIdx)
    $p2 = $d.Paragraphs.Item($secondIdx)

    # Whole range, including both paragraph marks, so the replacement XML's
    # own <w:p> elements fully take over (no stray proofErr/bookmark survives).
    $full = $d.Range($p1.Range.Start, $p2.Range.End)

    $newXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='00046CEC' w:rsidRDefault='00046CEC' w:rsidP='00046CEC'>" +
                "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
                "<w:r><w:t xml:space='preserve'>Create the </w:t></w:r>" +
                "<w:r><w:t>Object abstract class.</w:t></w:r>" +
              "</w:p>" +
              "<w:p>" +
                "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
                "<w:r><w:t>Create the Jewel</w:t></w:r>" +
                "<w:r><w:t xml:space='preserve'> class and inherit from the </w:t></w:r>" +
                "<w:r><w:t>Object class.</w:t></w:r>" +
              "</w:p>" +
              "<w:p>" +
                "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
                "<w:r><w:t>Create the Furniture class and inherit from the Object class.</w:t></w:r>" +
              "</w:p>" +
              "<w:p>" +
                "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
                "<w:r><w:t xml:space='preserve'>Create the </w:t></w:r>" +
                "<w:proofErr w:type='spellStart'/>" +
                "<w:r><w:t>GallifrianMirror</w:t></w:r>" +
                "<w:proofErr w:type='spellEnd'/>" +
                "<w:r><w:t xml:space='preserve'> class and inherit form the Object class.</w:t></w:r>" +
              "</w:p>" +
              "<w:p>" +
                "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
                "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
              "</w:p>" +
              "<w:p w:rsidR='00046CEC' w:rsidRDefault='00046CEC' w:rsidP='00046CEC'>" +
                "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
                "<w:proofErr w:type='spellStart'/>" +
                "<w:r><w:t>ff</w:t></w:r>" +
                "<w:proofErr w:type='spellEnd'/>" +
              "</w:p>"

    $full.InsertXML($newXml)

    Write-Host "Sprint 2 bullet list expanded. Paragraph count now:" $d.Paragraphs.Count
}
